# Update master to output generated at aa3dc9e
#
# Changes:
#  - Update the worksheet date heading from "2023-09-26 Tuesday" to
#    "2023-09-27 Wednesday".
#  - Replace all 25 division problems in the 5x5 grid with new values.
#
# The data rows of the table (1-indexed) are 1, 5, 9, 13, 17 (each
# followed by 3 blank spacer rows), with 5 columns each. Because several
# of the original problem strings repeat (and some new values coincide
# with other original values elsewhere in the document), we address each
# cell directly by its table position rather than relying on text
# search-and-replace, so there is no ambiguity about which occurrence is
# being updated.

$d = $word.ActiveDocument

# --- Heading date -----------------------------------------------------
$d.Content.Find.Execute("2023-09-26 Tuesday", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "2023-09-27 Wednesday", 2)

# --- Division grid ------------------------------------------------------
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "40÷6="
$t.Cell(1,2).Range.Text  = "38÷3="
$t.Cell(1,3).Range.Text  = "98÷9="
$t.Cell(1,4).Range.Text  = "68÷3="
$t.Cell(1,5).Range.Text  = "97÷8="

$t.Cell(5,1).Range.Text  = "15÷3="
$t.Cell(5,2).Range.Text  = "65÷5="
$t.Cell(5,3).Range.Text  = "21÷9="
$t.Cell(5,4).Range.Text  = "23÷3="
$t.Cell(5,5).Range.Text  = "54÷9="

$t.Cell(9,1).Range.Text  = "74÷3="
$t.Cell(9,2).Range.Text  = "98÷2="
$t.Cell(9,3).Range.Text  = "31÷2="
$t.Cell(9,4).Range.Text  = "94÷6="
$t.Cell(9,5).Range.Text  = "33÷2="

$t.Cell(13,1).Range.Text = "50÷2="
$t.Cell(13,2).Range.Text = "29÷5="
$t.Cell(13,3).Range.Text = "25÷2="
$t.Cell(13,4).Range.Text = "64÷9="
$t.Cell(13,5).Range.Text = "96÷5="

$t.Cell(17,1).Range.Text = "18÷6="
$t.Cell(17,2).Range.Text = "22÷7="
$t.Cell(17,3).Range.Text = "68÷5="
$t.Cell(17,4).Range.Text = "32÷4="
$t.Cell(17,5).Range.Text = "92÷2="
